# Reposition the Cassandra logo ("Picture 5") and the Apache Falcon logo
# ("Picture 27"), and remove the duplicate Apache Falcon logo ("Picture 4"),
# per the "add kibana, elasticsearch, change cassandra and zeppelin" commit.
#
# Shape.Left/.Top are expressed in points (1 pt = 12700 EMU) and are backed
# by single-precision floats that truncate toward zero when converted back
# to EMU on save, so nudge by half an EMU (in points) to land on the exact
# target EMU value.

$emuPerPt = 12700
$halfEmuInPt = 0.5 / $emuPerPt

function EmuToPt($emu) {
    return ($emu / $emuPerPt) + $halfEmuInPt
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Move "Picture 5" (Cassandra) from (1390116, 5671708) EMU to (555134, 5494462) EMU
$cassandra = $s.Shapes.Item("Picture 5")
$cassandra.Left = EmuToPt(555134)
$cassandra.Top = EmuToPt(5494462)

# Move "Picture 27" (Apache Falcon) from (696342, 4970555) EMU to (2459174, 5105944) EMU
$falcon27 = $s.Shapes.Item("Picture 27")
$falcon27.Left = EmuToPt(2459174)
$falcon27.Top = EmuToPt(5105944)

# Remove the extra duplicate "Picture 4" (Apache Falcon) shape entirely
$falcon4 = $s.Shapes.Item("Picture 4")
$falcon4.Delete()
